# Added Test Data For Hungary/Russia/Finland Market
#
# Adds three new country sheets (Russia, Finland, Hungary) after the existing
# "Denmark" sheet, cloning Denmark's layout/styles (same template used by all
# the other "Market" sheets in this workbook) and then swapping in the
# per-country NGC ticket code (B4) and market name (B2).

$wb = $excel.ActiveWorkbook
$denmark = $wb.Worksheets.Item("Denmark")

# --- Russia -----------------------------------------------------------
$denmark.Copy([System.Reflection.Missing]::Value, $denmark) | Out-Null
$russia = $wb.Worksheets.Item("Denmark (2)")
$russia.Name = "Russia"
$russia.Range("B4").Value = "NGC-2929/T2925/"
$russia.Range("B2").Value = "Russia Market"
$russia.Range("A1:D11").Select() | Out-Null

# --- Finland ------------------------------------------------------------
$denmark.Copy([System.Reflection.Missing]::Value, $russia) | Out-Null
$finland = $wb.Worksheets.Item("Denmark (2)")
$finland.Name = "Finland"
$finland.Range("B4").Value = "NGC-3130/T2957/"
$finland.Range("B2").Value = "Finland Market"
$finland.Range("A1:D11").Select() | Out-Null

# --- Hungary --------------------------------------------------------------
$denmark.Copy([System.Reflection.Missing]::Value, $finland) | Out-Null
$hungary = $wb.Worksheets.Item("Denmark (2)")
$hungary.Name = "Hungary"
$hungary.Range("B4").Value = "NGC-3104/T3004/"
$hungary.Range("B2").Value = "Hungary Market"

# Hungary becomes the new right-most / active tab.
$hungary.Range("E14").Select() | Out-Null
$hungary.Activate()
